# change leadscrew to 2mm, fix button BB range
$wb = $excel.ActiveWorkbook

$feed = $wb.Worksheets.Item("feed")

# Leadscrew pitch C6: 1.5 -> 2 (change leadscrew to 2mm)
$feed.Range("C6").Value = 2

# Fix button BB range A20: 0.06 -> 0.08
$feed.Range("A20").Value = 0.08

# Update selection on the feed sheet (no longer has a frozen topLeftCell, new active cell D9)
$feed.Activate()
$feed.Range("D9").Select()
